$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaText = $metaPara.Range.Text
if ($metaText -like "Meta description*") {
    $metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
    $metaRange.Delete()
}

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Beetle Mania Deluxe for Free - Slot
#    Game Review" right before the final ("Prompt for DALLE: ...") paragraph.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$insertPos = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParaXml = '<w:p ' + $xmlNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Beetle Mania Deluxe for Free - Slot Game Review</w:t></w:r></w:p>'
$spacerParaXml = '<w:p ' + $xmlNs + '/>'
$insertPos.InsertXML($newParaXml + $spacerParaXml)

# InsertXML leaves a spare empty paragraph behind the inserted one (it is
# what forces the real paragraph break) - remove it again.
$spacerPara = $d.Paragraphs($n + 1)
$spacerRange = $d.Range($spacerPara.Range.Start, $spacerPara.Range.End)
$spacerRange.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the text of the final paragraph (the old "Prompt for DALLE: ..."
#    text) with the new meta-description copy, keeping its italic run intact.
# ---------------------------------------------------------------------------
$n2 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($n2)
$finalTextRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalTextRange.Text = "Learn more about the insect-themed slot game Beetle Mania Deluxe and play it for free. Features include Wilds, Scatters, free spins, and a Gamble feature."
